# Remove the "lemmalist-greek" dependency row (row 10) from the Acknowledgments sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# This engine does not re-anchor Hyperlink objects when rows are deleted, so
# capture nothing here -- we rebuild the hyperlinks afterwards from scratch.
# First wipe all existing hyperlinks (per-item Delete is not wired up; only the
# whole-collection Delete works), then delete the row itself.
$ws.Hyperlinks.Delete()
$ws.Rows.Item(10).Delete()

# Recreate every remaining hyperlink at its (possibly shifted) new location.
$hl = $ws.Hyperlinks
$hl.Add($ws.Range("B2"), 'https://www.crummy.com/software/BeautifulSoup/')
$hl.Add($ws.Range("B4"), 'https://github.com/Ousret/charset_normalizer')
$hl.Add($ws.Range("B8"), 'https://github.com/Mimino666/langdetect')
$hl.Add($ws.Range("B9"), 'https://github.com/saffsd/langid.py')
$hl.Add($ws.Range("B11"), 'https://lxml.de/')
$hl.Add($ws.Range("B16"), 'https://www.numpy.org/')
$hl.Add($ws.Range("B18"), 'https://foss.heptapod.net/openpyxl/openpyxl')
$hl.Add($ws.Range("B26"), 'https://github.com/python-openxml/python-docx')
$hl.Add($ws.Range("B28"), 'https://github.com/psf/requests')
$hl.Add($ws.Range("B30"), 'https://scipy.org/scipylib/')
$hl.Add($ws.Range("F2"), 'https://bazaar.launchpad.net/~leonardr/beautifulsoup/bs4/view/head:/LICENSE')
$hl.Add($ws.Range("F4"), 'https://github.com/Ousret/charset_normalizer/blob/master/LICENSE')
$hl.Add($ws.Range("F8"), 'https://github.com/Mimino666/langdetect/blob/master/LICENSE')
$hl.Add($ws.Range("F9"), 'https://github.com/saffsd/langid.py/blob/master/LICENSE')
$hl.Add($ws.Range("F11"), 'https://github.com/lxml/lxml/blob/master/doc/licenses/BSD.txt')
$hl.Add($ws.Range("F16"), 'https://github.com/numpy/numpy/blob/master/LICENSE.txt')
$hl.Add($ws.Range("F18"), 'https://foss.heptapod.net/openpyxl/openpyxl/-/blob/branch/3.0/LICENCE.rst')
$hl.Add($ws.Range("F26"), 'https://github.com/python-openxml/python-docx/blob/master/LICENSE')
$hl.Add($ws.Range("F28"), 'https://github.com/requests/requests/blob/master/LICENSE')
$hl.Add($ws.Range("F30"), 'https://github.com/scipy/scipy/blob/master/LICENSE.txt')
$hl.Add($ws.Range("F25"), 'https://docs.python.org/3.8/license.html', 'psf-license-agreement-for-python-release')
$hl.Add($ws.Range("F20"), 'https://github.com/pyinstaller/pyinstaller/blob/develop/COPYING.txt')
$hl.Add($ws.Range("F23"), 'https://www.riverbankcomputing.com/static/Docs/PyQt5/introduction.html', 'license')
$hl.Add($ws.Range("B25"), 'https://www.python.org/')
$hl.Add($ws.Range("B20"), 'http://www.pyinstaller.org/')
$hl.Add($ws.Range("B23"), 'https://riverbankcomputing.com/software/pyqt/')
$hl.Add($ws.Range("B12"), 'https://matplotlib.org/')
$hl.Add($ws.Range("B14"), 'https://networkx.org/')
$hl.Add($ws.Range("F12"), 'https://matplotlib.org/users/license.html')
$hl.Add($ws.Range("F14"), 'https://github.com/networkx/networkx/blob/master/LICENSE.txt')
$hl.Add($ws.Range("B38"), 'https://github.com/amueller/word_cloud')
$hl.Add($ws.Range("F38"), 'https://github.com/amueller/word_cloud/blob/master/LICENSE')
$hl.Add($ws.Range("B3"), 'https://github.com/Esukhia/botok')
$hl.Add($ws.Range("B5"), 'https://github.com/cltk/cltk')
$hl.Add($ws.Range("B7"), 'https://github.com/fxsjy/jieba')
$hl.Add($ws.Range("B13"), 'https://github.com/taishi-i/nagisa')
$hl.Add($ws.Range("B15"), 'http://www.nltk.org/')
$hl.Add($ws.Range("B17"), 'https://github.com/yichen0831/opencc-python')
$hl.Add($ws.Range("B19"), 'https://github.com/lancopku/pkuseg-python')
$hl.Add($ws.Range("B22"), 'https://pyphen.org/')
$hl.Add($ws.Range("B21"), 'https://github.com/kmike/pymorphy2')
$hl.Add($ws.Range("B24"), 'https://github.com/PyThaiNLP/pythainlp')
$hl.Add($ws.Range("B27"), 'https://github.com/natasha/razdel')
$hl.Add($ws.Range("B29"), 'https://github.com/alvations/sacremoses')
$hl.Add($ws.Range("B31"), 'https://spacy.io/')
$hl.Add($ws.Range("B32"), 'https://github.com/ponrawee/ssg')
$hl.Add($ws.Range("B34"), 'https://github.com/fnl/syntok')
$hl.Add($ws.Range("B35"), 'https://github.com/sloria/TextBlob')
$hl.Add($ws.Range("B36"), 'https://github.com/mideind/Tokenizer')
$hl.Add($ws.Range("B37"), 'https://github.com/undertheseanlp/underthesea')
$hl.Add($ws.Range("F3"), 'https://github.com/Esukhia/botok/blob/master/LICENSE')
$hl.Add($ws.Range("F5"), 'https://github.com/cltk/cltk/blob/master/LICENSE')
$hl.Add($ws.Range("F7"), 'https://github.com/fxsjy/jieba/blob/master/LICENSE')
$hl.Add($ws.Range("F13"), 'https://github.com/taishi-i/nagisa/blob/master/LICENSE.txt')
$hl.Add($ws.Range("F15"), 'https://github.com/nltk/nltk/blob/develop/LICENSE.txt')
$hl.Add($ws.Range("F17"), 'https://github.com/yichen0831/opencc-python/blob/master/LICENSE.txt')
$hl.Add($ws.Range("F19"), 'https://github.com/lancopku/pkuseg-python/blob/master/LICENSE')
$hl.Add($ws.Range("F22"), 'https://github.com/Kozea/Pyphen/blob/master/LICENSE')
$hl.Add($ws.Range("F21"), 'https://github.com/kmike/pymorphy2/', 'pymorphy2')
$hl.Add($ws.Range("F24"), 'https://github.com/PyThaiNLP/pythainlp/blob/dev/LICENSE')
$hl.Add($ws.Range("F27"), 'https://github.com/natasha/razdel/blob/master/LICENSE')
$hl.Add($ws.Range("F29"), 'https://github.com/alvations/sacremoses/blob/master/LICENSE')
$hl.Add($ws.Range("F31"), 'https://github.com/explosion/spaCy/blob/master/LICENSE')
$hl.Add($ws.Range("F32"), 'https://github.com/ponrawee/ssg/blob/master/LICENSE')
$hl.Add($ws.Range("F34"), 'https://github.com/fnl/syntok/blob/master/LICENSE')
$hl.Add($ws.Range("F35"), 'https://github.com/sloria/TextBlob/blob/dev/LICENSE')
$hl.Add($ws.Range("F36"), 'https://github.com/mideind/Tokenizer/blob/master/LICENSE')
$hl.Add($ws.Range("F37"), 'https://github.com/undertheseanlp/underthesea/blob/master/LICENSE')
$hl.Add($ws.Range("B6"), 'https://github.com/Xangis/extra-stopwords')
$hl.Add($ws.Range("B10"), 'https://github.com/michmech/lemmatization-lists')
$hl.Add($ws.Range("B33"), 'https://github.com/stopwords-iso/stopwords-iso')
$hl.Add($ws.Range("F6"), 'https://github.com/Xangis/extra-stopwords/blob/master/LICENSE')
$hl.Add($ws.Range("F10"), 'https://github.com/michmech/lemmatization-lists/blob/master/LICENCE')
$hl.Add($ws.Range("F33"), 'https://github.com/stopwords-iso/stopwords-iso/blob/master/LICENSE')

# Restore the view roughly where the user left it after the edit.
$ws.Range("D11").Select()

